$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Atualizei dados bibi - linha de 2025 (row 9)
$ws.Range("B9").Value = 3999051.08
$ws.Range("C9").Value = 625176.24
$ws.Range("D9").Value = 4624227.32
$ws.Range("E9").Value = 13.51958276999237
$ws.Range("F9").Value = 86.48041723000762
$ws.Range("G9").Value = -39.57983367644407
$ws.Range("H9").Value = -27.78266514202331
$ws.Range("I9").Value = 40123
$ws.Range("J9").Value = 1718
$ws.Range("K9").Value = 41841
$ws.Range("L9").Value = 28990
$ws.Range("M9").Value = 159.5111183166609
$ws.Range("N9").Value = 8.90148935408015
